$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.813.72"
$ws.Range("E2").Value = "  -3.47%  "
$ws.Range("D3").Value = "2.482.52"
$ws.Range("E3").Value = "  -6.16%  "
$ws.Range("E4").Value = "  +0.06%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "556.86"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -4.22%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "148.14"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -5.43%  "
$ws.Range("E7").Value = "  +0.05%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.602"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").Value = "2.480.77"
$ws.Range("E9").Value = "  -6.17%  "
$ws.Range("E10").Value = "  -8.23%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "5.51"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -5.43%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.360"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -6.25%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "26.65"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -6.84%  "
$ws.Range("D15").Value = "2.934.87"
$ws.Range("E15").Value = "  -6.04%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.0000169"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -8.72%  "
$ws.Range("D17").Value = "61.762.13"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "2.484.82"
$ws.Range("E18").Value = "  -6.05%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "11.27"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -7.77%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "7.21"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -6.96%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "4.25"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -6.24%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "323.42"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "1.89"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +2.26%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "64.56"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -5.52%  "
$ws.Range("E26").Value = "  -9.84%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "571.22"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("D28").Value = "2.607.19"
$ws.Range("E28").Value = "  -6.15%  "
$ws.Range("E29").Value = "  -7.46%  "
$ws.Range("E30").Value = "  -0.02%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "8.39"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -10.38%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "7.82"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -4.58%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "0.151"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -6.43%  "
$ws.Range("E34").Value = "  -6.43%  "
$ws.Range("E35").Value = "  -7.94%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "5.99"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -9.87%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "4.97"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -9.26%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.386"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -4.53%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "18.67"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -5.68%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "144.77"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -4.61%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "1.78"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -7.17%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "2.46"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -4.56%  "
$ws.Range("E45").Value = "  -2.84%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "149.42"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -8.80%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "3.67"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -6.35%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "22.12"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -9.78%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.0545"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -8.03%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.600"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -5.63%  "
$ws.Range("E51").Value = "  -5.40%  "
